$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the USD Amount value in T2 (246465 -> 248216)
$ws.Range("T2").Value = 248216

# Move / leave the active selection on T2 (was T3 in the original file)
[void]$ws.Range("T2").Select()
